$d = $word.ActiveDocument
$d.Content.Find.Execute("As a whole minority groups", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The Bakken Oil Boom started around within the state of North Dakota and Montana, as well as parts of Canada. Native Americans", 2)
